$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 433, shifting existing rows 433..510 down to 434..511
$ws.Rows.Item(433).Insert()

# Populate the new row 433 with a copy of the (now shifted) row 434 data, then apply the
# specific field changes from the commit.
$ws.Range("A433").Value2 = 10
$ws.Range("B433").Value2 = "Vega Modelo de Temuco"
$ws.Range("C433").Value2 = "La Araucanía"
$ws.Range("D433").Value2 = 45180
$ws.Range("D433").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E433").Value2 = 9
$ws.Range("F433").Value2 = 100112001
$ws.Range("G433").Value2 = "Berenjena"
$ws.Range("H433").Value2 = "Sin especificar"
$ws.Range("I433").Value2 = "Primera"
$ws.Range("J433").Value2 = 200
$ws.Range("K433").Value2 = 10000
$ws.Range("L433").Value2 = 12000
$ws.Range("M433").Value2 = 11200
$ws.Range("N433").Value2 = "$/caja 40 unidades"
$ws.Range("O433").Value2 = "Región de Arica y Parinacota"
$ws.Range("P433").Value2 = 280
$ws.Range("Q433").Value2 = 40
$ws.Range("R433").Value2 = "Hortaliza"
